# -----------------------------------------------------------------------
# Refresh the cryptocurrency price/volume snapshot on Sheet1.
#
# Each row holds one coin's data scraped from coinranking.com:
#   column D = current price, column E = 1-hour volume change (%).
# This run updates the rows whose price and/or 1h-change figures moved
# since the previous GitHub Actions scrape.
#
# Price strings that look like plain numbers (e.g. "585.92") are written
# with a leading apostrophe so Excel stores them as text, exactly like
# the existing cells, instead of silently converting them to numbers.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Bitcoin
$ws.Range("D2").Value = '63.300.62'
$ws.Range("E2").Value = '  +6.32%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '3.122.56'
$ws.Range("E3").Value = '  +4.26%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.03%  '

# Row 5: BNB
$ws.Range("D5").Value = '''585.92'
$ws.Range("E5").Value = '  +4.06%  '

# Row 6: Solana
$ws.Range("D6").Value = '''145.03'
$ws.Range("E6").Value = '  +4.28%  '

# Row 8: LidoStakedEther
$ws.Range("D8").Value = '3.113.64'
$ws.Range("E8").Value = '  +4.31%  '

# Row 9: XRP
$ws.Range("E9").Value = '  +1.85%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  +12.63%  '

# Row 11: Toncoin
$ws.Range("D11").Value = '''5.79'
$ws.Range("E11").Value = '  +9.77%  '

# Row 12: Cardano
$ws.Range("E12").Value = '  +3.15%  '

# Row 13: ShibaInu
$ws.Range("D13").Value = '''0.0000248'
$ws.Range("E13").Value = '  +7.93%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''35.64'
$ws.Range("E14").Value = '  +5.16%  '

# Row 15: TRON
$ws.Range("E15").Value = '  +0.42%  '

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '3.637.62'
$ws.Range("E16").Value = '  +4.16%  '

# Row 17: Polkadot
$ws.Range("E17").Value = '  -0.16%  '

# Row 18: WrappedBTC
$ws.Range("D18").Value = '63.210.12'
$ws.Range("E18").Value = '  +6.21%  '

# Row 19: WrappedEther
$ws.Range("D19").Value = '3.119.00'
$ws.Range("E19").Value = '  +4.28%  '

# Row 20: BitcoinCash
$ws.Range("D20").Value = '''467.51'

# Row 21: Chainlink
$ws.Range("E21").Value = '  +3.62%  '

# Row 22: Polygon
$ws.Range("E22").Value = '  +1.03%  '

# Row 23: Uniswap
$ws.Range("D23").Value = '''7.55'
$ws.Range("E23").Value = '  +6.49%  '

# Row 24: InternetComputer(DFINITY)
$ws.Range("D24").Value = '''13.31'
$ws.Range("E24").Value = '  -1.74%  '

# Row 25: Litecoin
$ws.Range("E25").Value = '  +2.25%  '

# Row 26: Dai
$ws.Range("E26").Value = '  +0.04%  '

# Row 27: ImmutableX
$ws.Range("E27").Value = '  +0.51%  '

# Row 28: RenderToken
$ws.Range("E28").Value = '  +7.31%  '

# Row 29: PancakeSwap
$ws.Range("E29").Value = '  +5.37%  '

# Row 30: FirstDigitalUSD
$ws.Range("E30").Value = '  +0.02%  '

# Row 31: NEARProtocol
$ws.Range("D31").Value = '''6.86'
$ws.Range("E31").Value = '  +9.69%  '

# Row 32: EthereumClassic
$ws.Range("D32").Value = '''26.98'
$ws.Range("E32").Value = '  +4.65%  '

# Row 33: Hedera
$ws.Range("E33").Value = '  +3.18%  '

# Row 34: PEPE
$ws.Range("D34").Value = '0.0₃0871'
$ws.Range("E34").Value = '  +11.85%  '

# Row 35: Stacks
$ws.Range("E35").Value = '  +16.41%  '

# Row 36: Mantle
$ws.Range("E36").Value = '  +5.26%  '

# Row 37: dogwifhat
$ws.Range("D37").Value = '''3.32'
$ws.Range("E37").Value = '  +18.86%  '

# Row 38: Filecoin
$ws.Range("E38").Value = '  +2.26%  '

# Row 39: OKB
$ws.Range("E39").Value = '  +3.98%  '

# Row 40: Bittensor
$ws.Range("D40").Value = '''432.56'
$ws.Range("E40").Value = '  +7.77%  '

# Row 41: Cosmos
$ws.Range("E41").Value = '  +1.78%  '

# Row 42: Maker
$ws.Range("D42").Value = '2.930.62'
$ws.Range("E42").Value = '  +6.02%  '

# Row 43: VeChain
$ws.Range("E43").Value = '  +4.49%  '

# Row 44: TheGraph
$ws.Range("E44").Value = '  +10.97%  '

# Row 45: Kaspa
$ws.Range("E45").Value = '  +5.70%  '

# Row 46: Fetch.AI
$ws.Range("E46").Value = '  +7.57%  '

# Row 47: Arweave
$ws.Range("D47").Value = '''35.34'
$ws.Range("E47").Value = '  +2.79%  '

# Row 48: USDe
$ws.Range("E48").Value = '  -0.01%  '

# Row 49: Monero
$ws.Range("D49").Value = '''123.41'
$ws.Range("E49").Value = '  +0.36%  '

# Row 50: Stellar
$ws.Range("E50").Value = '  +0.90%  '

# Row 51: InjectiveProtocol
$ws.Range("E51").Value = '  +4.06%  '

